$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => [D_new, E_new] (D left $null when unchanged per the source diff)
$updates = @(
    @{ Row = 2;  D = "28.196.38";  E = "  +3.58%  " },
    @{ Row = 3;  D = "1.577.58";   E = "  +0.41%  " },
    @{ Row = 4;  D = $null;        E = "  -1.01%  " },
    @{ Row = 5;  D = "213.23";     E = "  +0.76%  " },
    @{ Row = 6;  D = "0.494";      E = "  +0.17%  " },
    @{ Row = 7;  D = $null;        E = "  -0.95%  " },
    @{ Row = 8;  D = "23.53";      E = "  +6.28%  " },
    @{ Row = 9;  D = $null;        E = "  +0.66%  " },
    @{ Row = 11; D = "0.0886";     E = "  +1.95%  " },
    @{ Row = 12; D = "1.801.94";   E = "  +0.41%  " },
    @{ Row = 13; D = "1.567.24";   E = "  -0.21%  " },
    @{ Row = 14; D = $null;        E = "  -0.82%  " },
    @{ Row = 15; D = "0.526";      E = "  +1.12%  " },
    @{ Row = 16; D = "28.150.50";  E = "  +3.40%  " },
    @{ Row = 17; D = "63.78";      E = "  +2.16%  " },
    @{ Row = 18; D = "230.17";     E = "  +6.34%  " },
    @{ Row = 19; D = $null;        E = "  +0.43%  " },
    @{ Row = 20; D = "7.46";       E = "  +0.49%  " },
    @{ Row = 21; D = $null;        E = "  -0.99%  " },
    @{ Row = 22; D = $null;        E = "  -0.64%  " },
    @{ Row = 23; D = "9.33";       E = "  +0.89%  " },
    @{ Row = 24; D = $null;        E = "  -0.52%  " },
    @{ Row = 25; D = "152.24";     E = "  -1.17%  " },
    @{ Row = 26; D = "15.24";      E = "  +0.92%  " },
    @{ Row = 27; D = "6.59";       E = "  -2.20%  " },
    @{ Row = 28; D = "0.107";      E = "  -0.11%  " },
    @{ Row = 29; D = $null;        E = "  -0.90%  " },
    @{ Row = 30; D = $null;        E = "  +0.03%  " },
    @{ Row = 31; D = $null;        E = "  +0.02%  " },
    @{ Row = 32; D = $null;        E = "  -0.91%  " },
    @{ Row = 33; D = $null;        E = "  -0.93%  " },
    @{ Row = 34; D = "1.415.98";   E = "  -2.53%  " },
    @{ Row = 35; D = $null;        E = "  -1.36%  " },
    @{ Row = 36; D = $null;        E = "  -4.89%  " },
    @{ Row = 37; D = "2.32";       E = "  -1.33%  " },
    @{ Row = 38; D = $null;        E = "  -0.38%  " },
    @{ Row = 39; D = "2.53";       E = "  +7.55%  " },
    @{ Row = 40; D = $null;        E = "  +0.96%  " },
    @{ Row = 41; D = $null;        E = "  -0.22%  " },
    @{ Row = 42; D = $null;        E = "  -1.07%  " },
    @{ Row = 43; D = $null;        E = "  -3.36%  " },
    @{ Row = 44; D = "0.973";      E = "  -2.64%  " },
    @{ Row = 45; D = $null;        E = "  +4.52%  " },
    @{ Row = 46; D = "63.81";      E = "  -1.45%  " },
    @{ Row = 47; D = "1.714.06";   E = "  +0.44%  " },
    @{ Row = 48; D = "87.10";      E = "  +1.32%  " },
    @{ Row = 49; D = "0.0₆0107";   E = "  +2.38%  " },
    @{ Row = 50; D = "0.0526";     E = "  +0.92%  " },
    @{ Row = 51; D = "0.0944";     E = "  -1.53%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Leading apostrophe forces these price strings (which can look like
        # numbers/dates, e.g. "213.23") to be stored as text, matching the
        # original inlineStr cell type; resetting the style afterwards keeps
        # the cell's formatting untouched (no quotePrefix styling left behind).
        $dcell = $ws.Cells.Item($u.Row, 4)
        $dcell.Value = "'" + $u.D
        $dcell.Style = "Normal"
    }
    $ecell = $ws.Cells.Item($u.Row, 5)
    $ecell.Value = $u.E
}
